# Update cryptos list values per the data diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values in column D are plain decimal numbers (e.g. "96.00").
# The source data stores these as TEXT (to preserve exact formatting, like
# trailing zeros), so we mark those specific cells as Text ("@") before
# assigning the values. This keeps Excel from auto-converting the strings
# into numeric doubles (which would corrupt values like "96.00" -> 96).
$textCells = @(
    'D5',
    'D6',
    'D7',
    'D8',
    'D9',
    'D10',
    'D11',
    'D12',
    'D14',
    'D17',
    'D19',
    'D21',
    'D22',
    'D23',
    'D24',
    'D28',
    'D29',
    'D31',
    'D32',
    'D33',
    'D34',
    'D35',
    'D36',
    'D38',
    'D44',
    'D45',
    'D46',
    'D47',
    'D48',
    'D50'
)
foreach ($cellAddr in $textCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

$ws.Range('D2').Value = '42.654.63'
$ws.Range('E2').Value = '  -1.00%  '

$ws.Range('D3').Value = '2.280.74'
$ws.Range('E3').Value = '  -1.05%  '

$ws.Range('E4').Value = '  -0.11%  '

$ws.Range('D5').Value = '304.81'
$ws.Range('E5').Value = '  +1.50%  '

$ws.Range('D6').Value = '96.00'
$ws.Range('E6').Value = '  -2.16%  '

$ws.Range('D7').Value = '0.507'
$ws.Range('E7').Value = '  -2.35%  '

$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('D9').Value = '0.500'
$ws.Range('E9').Value = '  -3.13%  '

$ws.Range('D10').Value = '35.26'
$ws.Range('E10').Value = '  -2.34%  '

$ws.Range('D11').Value = '0.0790'
$ws.Range('E11').Value = '  -0.18%  '

$ws.Range('D12').Value = '18.25'
$ws.Range('E12').Value = '  +3.21%  '

$ws.Range('E13').Value = '  +0.88%  '

$ws.Range('D14').Value = '6.74'
$ws.Range('E14').Value = '  -2.23%  '

$ws.Range('D15').Value = '2.631.87'
$ws.Range('E15').Value = '  -1.11%  '

$ws.Range('D16').Value = '2.262.51'
$ws.Range('E16').Value = '  -1.37%  '

$ws.Range('D17').Value = '0.777'
$ws.Range('E17').Value = '  -1.55%  '

$ws.Range('D18').Value = '42.577.22'
$ws.Range('E18').Value = '  -0.89%  '

$ws.Range('D19').Value = '12.93'
$ws.Range('E19').Value = '  +1.44%  '

$ws.Range('D20').Value = '0.0₃0893'
$ws.Range('E20').Value = '  -2.17%  '

$ws.Range('D21').Value = '6.02'
$ws.Range('E21').Value = '  -2.27%  '

$ws.Range('D22').Value = '67.13'
$ws.Range('E22').Value = '  -1.86%  '

$ws.Range('D23').Value = '235.64'
$ws.Range('E23').Value = '  -0.98%  '

$ws.Range('D24').Value = '2.13'
$ws.Range('E24').Value = '  -1.54%  '

$ws.Range('E25').Value = '  +0.16%  '

$ws.Range('E26').Value = '  +0.59%  '

$ws.Range('E27').Value = '  -0.07%  '

$ws.Range('D28').Value = '25.15'
$ws.Range('E28').Value = '  +0.20%  '

$ws.Range('D29').Value = '166.32'
$ws.Range('E29').Value = '  +1.28%  '

$ws.Range('E30').Value = '  +0.59%  '

$ws.Range('D31').Value = '9.04'
$ws.Range('E31').Value = '  -1.15%  '

$ws.Range('D32').Value = '32.94'
$ws.Range('E32').Value = '  -0.63%  '

$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  -0.05%  '

$ws.Range('D34').Value = '4.75'
$ws.Range('E34').Value = '  -1.35%  '

$ws.Range('D35').Value = '4.96'
$ws.Range('E35').Value = '  -3.20%  '

$ws.Range('D36').Value = '17.55'
$ws.Range('E36').Value = '  -3.01%  '

$ws.Range('E37').Value = '  -1.22%  '

$ws.Range('D38').Value = '0.0690'
$ws.Range('E38').Value = '  -1.20%  '

$ws.Range('E39').Value = '  -0.87%  '

$ws.Range('E40').Value = '  -2.17%  '

$ws.Range('E41').Value = '  -1.54%  '

$ws.Range('E42').Value = '  -3.40%  '

$ws.Range('D43').Value = '2.006.44'
$ws.Range('E43').Value = '  -0.71%  '

$ws.Range('D44').Value = '0.0279'
$ws.Range('E44').Value = '  -2.85%  '

$ws.Range('D45').Value = '18.02'
$ws.Range('E45').Value = '  +3.10%  '

$ws.Range('D46').Value = '9.98'
$ws.Range('E46').Value = '  -3.92%  '

$ws.Range('D47').Value = '2.08'
$ws.Range('E47').Value = '  -6.36%  '

$ws.Range('D48').Value = '2.77'
$ws.Range('E48').Value = '  -2.82%  '

$ws.Range('E49').Value = '  +5.71%  '

$ws.Range('D50').Value = '53.55'
$ws.Range('E50').Value = '  -1.51%  '

$ws.Range('D51').Value = '2.500.23'
$ws.Range('E51').Value = '  -1.00%  '

